$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Mrs. N. Kasirathi
$ws.Range("J2").Value = "VEC-007-04-187"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.researchgate.net/profile/Kasirathi-Natarajan?ev=hdr_xprf")

# Row 3 - Dr. K. Vidhya
$ws.Range("J3").Value = "VEC-007-04-178"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.researchgate.net/profile/Vidhya-Kumar-4?ev=hdr_xprf ")

# Row 8 - Dr. R. Geetha
$ws.Range("J8").Value = "VEC-007-04-185"

# Row 13 - Mr. V. Jagadeesh Babu
$ws.Range("J13").Value = "VEC-007-04-184"

# Row 4 - Dr. D. Sunitha
$ws.Range("E4").Value = "https://www.researchgate.net/scientific-contributions/D-Sunitha-2096230154"

# Row 5 - Mrs. G. Angalaparameswari
$ws.Range("E5").Value = "https://www.researchgate.net/profile/Angala-Gandhinathan"

# Row 7 - Dr. S. Raghuraman
$ws.Range("E7").Value = "https://www.researchgate.net/profile/Raghuraman-Sivalingam"

# Column width tweaks (column 11 is intentionally left untouched so it naturally
# splits off from column 10 while keeping its original width of 19)
$ws.Columns.Item(4).ColumnWidth = 56.250354
$ws.Columns.Item(5).ColumnWidth = 59.25051
$ws.Columns.Item(6).ColumnWidth = 30.751479
$ws.Columns.Item(7).ColumnWidth = 19.583792
$ws.Columns.Item(9).ColumnWidth = 45.250354
$ws.Columns.Item(10).ColumnWidth = 45.916917

# Selection update
$ws.Range("A10:B15").Select()
